$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1150.01174265776
$ws.Range("I2").Value = 466.011742657762

$ws.Range("B3").Value = 794.570596431729
$ws.Range("I3").Value = 420.570596431729

$ws.Range("B4").Value = 746.459753041819
$ws.Range("I4").Value = 313.459753041819

$ws.Range("B5").Value = 709.301964513145
$ws.Range("I5").Value = 169.301964513145

$ws.Range("B6").Value = 858.13586424188
$ws.Range("I6").Value = 172.13586424188

$ws.Range("B7").Value = 944.922434812283
$ws.Range("I7").Value = 118.922434812283

$ws.Range("B8").Value = 754.391524325922
$ws.Range("I8").Value = 441.391524325922

$ws.Range("B9").Value = 467.559927901096
$ws.Range("I9").Value = 154.559927901096

$ws.Range("B10").Value = 367.115758852823
$ws.Range("I10").Value = 47.1157588528231

$ws.Range("B11").Value = 614.723652321585
$ws.Range("I11").Value = 3.72365232158518
$ws.Range("J11").Value = "Decrease"

$ws.Range("B12").Value = 1274.49620447927
$ws.Range("I12").Value = -521.503795520731

$ws.Range("B13").Value = 1321.62745765955
$ws.Range("I13").Value = -138.37254234045

$ws.Range("B14").Value = 1172.46864860675
$ws.Range("I14").Value = 610.468648606752

$ws.Range("B15").Value = 864.458535230527
$ws.Range("I15").Value = 499.458535230527

$ws.Range("B16").Value = 768.268637250961
$ws.Range("I16").Value = 344.268637250961

$ws.Range("B17").Value = 720.226293035214
$ws.Range("I17").Value = 210.226293035214

$ws.Range("B18").Value = 809.099196005296
$ws.Range("I18").Value = 124.099196005296

$ws.Range("B19").Value = 888.90498475726
$ws.Range("I19").Value = 164.90498475726

$ws.Range("B20").Value = 778.917171975433
$ws.Range("I20").Value = 311.917171975433

$ws.Range("B21").Value = 522.022828057113
$ws.Range("I21").Value = 267.022828057113

$ws.Range("B22").Value = 390.352079936385
$ws.Range("I22").Value = 99.3520799363846

$ws.Range("B23").Value = 588.670402469447
$ws.Range("I23").Value = -89.3295975305531

$ws.Range("B24").Value = 1141.77430319651
$ws.Range("I24").Value = -978.225696803489

$ws.Range("B25").Value = 1256.79106816691
$ws.Range("I25").Value = -1145.20893183309

$ws.Range("B26").Value = 1174.40575216711
$ws.Range("I26").Value = 393.405752167114

$ws.Range("B27").Value = 918.821945822292
$ws.Range("I27").Value = 588.821945822292

$ws.Range("B28").Value = 794.455881855354
$ws.Range("I28").Value = 450.455881855354

$ws.Range("B29").Value = 733.924072834326
$ws.Range("I29").Value = 400.924072834326

$ws.Range("B30").Value = 780.649259438941
$ws.Range("I30").Value = 227.649259438941

$ws.Range("B31").Value = 843.721702418014
$ws.Range("I31").Value = 277.721702418014

$ws.Range("B32").Value = 785.797295295088
$ws.Range("I32").Value = 381.797295295088

$ws.Range("B33").Value = 572.438308811939
$ws.Range("I33").Value = 337.438308811939

$ws.Range("B34").Value = 415.92931577723
$ws.Range("I34").Value = 210.92931577723

$ws.Range("B35").Value = 532.170424641257
$ws.Range("I35").Value = 132.170424641257

$ws.Range("B36").Value = 1025.17076006902
$ws.Range("I36").Value = 130.170760069023

$ws.Range("B37").Value = 1177.64198979482
$ws.Range("I37").Value = 665.641989794823

$ws.Range("B38").Value = 1156.35835140193
$ws.Range("I38").Value = 939.358351401931

$ws.Range("B39").Value = 955.684846570057
$ws.Range("I39").Value = 685.684846570057

$ws.Range("B40").Value = 819.60968732024
$ws.Range("I40").Value = 489.60968732024

$ws.Range("B41").Value = 749.34564656206
$ws.Range("I41").Value = 444.34564656206

$ws.Range("B42").Value = 765.368027079837
$ws.Range("I42").Value = 366.368027079837

$ws.Range("B43").Value = 810.355332157052
$ws.Range("I43").Value = 445.355332157052

$ws.Range("B44").Value = 782.689736320927
$ws.Range("I44").Value = 438.689736320927

$ws.Range("B45").Value = 615.750687717448
$ws.Range("I45").Value = 375.750687717448

$ws.Range("B46").Value = 443.866694311528
$ws.Range("I46").Value = 245.866694311528

$ws.Range("B47").Value = 496.118654569495
$ws.Range("I47").Value = 57.1186545694946

$ws.Range("B48").Value = 934.259792349094
$ws.Range("I48").Value = -385.740207650906

$ws.Range("B49").Value = 1094.15076186458
$ws.Range("I49").Value = -27.8492381354176
